# Rewrite the "closest_residues_comb" data table (rows 2-22) on the active
# sheet with the draft-version values: 21 data rows under the existing
# header row (Residue Combination | Count | Permeated Residues | Frames).
#
# Columns A, C and D hold free-text (residue lists / frame numbers), many of
# which look like plain integers ("130", "780", "5552", ...). Force those two
# columns to text *before* writing so Excel doesn't silently reinterpret them
# as numbers, then drop back to the Normal style so no stray number format
# lingers on the cells once the text is committed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('98, 130, 748, 1073', 3, '130, 130, 130', '5178, 5582, 6488'),
    @('98, 455, 780, 1105', 1, '780', '5552'),
    @('130, 780, 780', 1, '130', '3171'),
    @('98, 130, 748, 780', 1, '130', '3631'),
    @('98, 780, 1105', 1, '780', '3666'),
    @('130, 455, 748, 1073', 1, '130', '5131'),
    @('98, 130, 1073, SF', 1, '130', '4415'),
    @('423, 1073, 1105, SF, SF', 1, '1105', '4994'),
    @('130, 455, 748, 780', 1, '130', '5269'),
    @('130, 748, 1073, 1105', 1, '1105', '5399'),
    @('423, 748, 780, 1073', 1, '780', '5677'),
    @('423, 748, 1073, 1105', 2, '1105, 1105', '5331, 5433'),
    @('130, 423, 748, 1073', 1, '130', '5886'),
    @('98, 130, 455, 748, 1073', 1, '130', '6016'),
    @('130, 423, 748, 1073, SF', 1, '130', '6202'),
    @('130, 130, 423', 1, '130', '6426'),
    @('130, 455, 780, 1105', 1, '130', '6561'),
    @('98, 780, 780, 1105', 1, '780', '6359'),
    @('98, 130, 130, 455, 780', 1, '130', '6727'),
    @('130, 130, 423, 1073', 1, '130', '6670'),
    @('98, 98, 455, 455', 1, '455', '6748')
)

$rowCount = $data.Count
$lastRow = 1 + $rowCount

# Force columns A, C, D to text so numeric-looking strings (e.g. "780",
# "5552") round-trip as text instead of being coerced to numbers. (Two
# separate calls -- a single comma-joined multi-area range doesn't reliably
# apply the format to every area.)
$colARange = $ws.Range("A2:A$lastRow")
$colCDRange = $ws.Range("C2:D$lastRow")
$colARange.NumberFormat = "@"
$colCDRange.NumberFormat = "@"

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

# Drop the temporary text format back to Normal now that the literal text is
# committed, so no stray number-format artifact is left on the cells.
$colARange.Style = "Normal"
$colCDRange.Style = "Normal"
